$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01514828764759746, 0.04240448674262143, 26.21740644021617, 2367095152636972, 2367095152636998)
    3 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 645.3272768299601, 646.7503648927143)
    4 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 8.660232485948974, 14.90378790461981)
    5 = @(1.459612070389937, 10.29869402782916, 0.1575252929769615, 616238.5361209477, 616250.4519523389)
    6 = @(0.003994804209775715, 0.002777888934908601, 26.21740644021617, 645.3272768299601, 671.5514559633209)
    7 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 11.80884686099532)
    8 = @(1.459612070389937, 1.667794583268128, 337.1190423067083, 616238.5361209477, 616578.7825699081)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
}
